$wb = $excel.ActiveWorkbook

$oldGuid = "940c231c-cdc6-458a-8955-ab89043e0388"
$newGuid = "a773750d-5825-4da6-af91-1b3ae8aa421e"

$oldZhCnHash = "989096b01a9e7fe8fde66d198afa168a9d196fcb"
$newZhCnHash = "b8803bb11e41f98f94d5c944019c8932b865a3bd"
$oldDeDeHash = "989096b01a9e7fe8fde66d198afa168a9d196fcb"
$newDeDeHash = "b8803bb11e41f98f94d5c944019c8932b865a3bd"

# ---------------------------------------------------------------
# Overview sheet: new source file name/path + refreshed HO date
# ---------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("A2").Value = "$newGuid.md"
$wsOverview.Range("B2").Value = "e2e\$newGuid.md"
foreach ($hl in $wsOverview.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$B$2') {
        $hl.TextToDisplay = "e2e\$newGuid.md"
    }
}
$wsOverview.Range("G2").Value = "2016-08-14 17:19:57"

# ---------------------------------------------------------------
# zh-cn sheet: new handoff round, handback not received yet
# ---------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("A2").Value = "$newGuid.md"
foreach ($hl in $wsZhCn.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$A$2') {
        $hl.TextToDisplay = "$newGuid.md"
    }
}
$wsZhCn.Range("G2").Value = "$newGuid.$newZhCnHash.zh-cn.xlf"
$wsZhCn.Range("H2").Value = "2016-08-14 17:19:49"

foreach ($hl in $wsZhCn.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$I$2') {
        $hl.Delete()
    }
}
$wsZhCn.Range("I2").Style = "Normal"
$wsZhCn.Range("I2").Value = ""
$wsZhCn.Range("J2").Value = ""
$wsZhCn.Range("K2").Value = "0001-01-01 00:00:00"

# ---------------------------------------------------------------
# de-de sheet: new handoff round, handback not received yet
# ---------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("A2").Value = "$newGuid.md"
foreach ($hl in $wsDeDe.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$A$2') {
        $hl.TextToDisplay = "$newGuid.md"
    }
}
$wsDeDe.Range("G2").Value = "$newGuid.$newDeDeHash.de-de.xlf"
$wsDeDe.Range("H2").Value = "2016-08-14 17:19:57"

foreach ($hl in $wsDeDe.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$I$2') {
        $hl.Delete()
    }
}
$wsDeDe.Range("I2").Style = "Normal"
$wsDeDe.Range("I2").Value = ""
$wsDeDe.Range("J2").Value = ""
$wsDeDe.Range("K2").Value = "0001-01-01 00:00:00"

# ---------------------------------------------------------------
# Column width tweaks on zh-cn / de-de (cols I & J no longer need
# to be as wide now that the target/handback file columns are empty)
# ---------------------------------------------------------------
$wsZhCn.Columns.Item(9).ColumnWidth = 18.6506053379604
$wsZhCn.Columns.Item(10).ColumnWidth = 21.7054770333426
$wsDeDe.Columns.Item(9).ColumnWidth = 18.6506053379604
$wsDeDe.Columns.Item(10).ColumnWidth = 21.7054770333426
